$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Min values for the first three parameter rows ---
$ws.Range("B2").Value = 5       # alpha_distance_range Min: 3.8 -> 5
$ws.Range("B3").Value = 5.5     # beta_distance_range Min: 5.3 -> 5.5
$ws.Range("B4").Value = 0.7     # ratio_threshold_range Min: 0.8 -> 0.7

# --- Remove the "theta_threshold_range" row entirely ---
# This shifts the old row 6 ("pie_threshold_range") up to become row 5,
# carrying its own text/values/styles with it, and drops the dimension
# from A1:C6 down to A1:C5.
$ws.Rows("5").Delete()

# --- Normalize formatting on the (now shifted) pie_threshold_range row ---
# Previously the "125" cell (old C5) used a special Times New Roman style,
# and B6 used that same special style; after the shift make both B5/C5
# use the same plain style as the rest of the data rows.
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)

# --- Update the pie_threshold_range Min/Max values ---
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# --- Update the active selection to match the saved view state ---
[void]$ws.Range("C3").Select()
